# Auto-generated Excel COM-interop script to apply the diff to 上海-漫展信息.xlsx
# Sheets: 1=展览, 2=演出, 3=本地生活(本地生活 has a row removed + reshuffle), 4=全部类型

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value2 = $Text
    $Cell.Style = "Normal"
}

function Set-NumCell {
    param($Cell, $Number)
    $Cell.Value2 = $Number
}

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item(1)
Set-NumCell $ws.Cells.Item(3, 6) 1712
Set-NumCell $ws.Cells.Item(4, 6) 1148
Set-NumCell $ws.Cells.Item(6, 6) 142
Set-NumCell $ws.Cells.Item(7, 6) 1401
Set-NumCell $ws.Cells.Item(8, 6) 63
Set-NumCell $ws.Cells.Item(9, 6) 93
Set-NumCell $ws.Cells.Item(10, 6) 627
Set-NumCell $ws.Cells.Item(12, 6) 87
Set-NumCell $ws.Cells.Item(13, 6) 1310
Set-NumCell $ws.Cells.Item(14, 6) 462
Set-NumCell $ws.Cells.Item(15, 6) 465
Set-NumCell $ws.Cells.Item(16, 6) 127
Set-NumCell $ws.Cells.Item(17, 6) 25
Set-NumCell $ws.Cells.Item(18, 6) 685
Set-NumCell $ws.Cells.Item(19, 6) 2539
Set-NumCell $ws.Cells.Item(21, 6) 46
Set-NumCell $ws.Cells.Item(22, 6) 16
Set-NumCell $ws.Cells.Item(24, 6) 283
Set-NumCell $ws.Cells.Item(27, 6) 97
Set-NumCell $ws.Cells.Item(28, 6) 562
Set-NumCell $ws.Cells.Item(29, 6) 918
Set-NumCell $ws.Cells.Item(31, 6) 63
Set-NumCell $ws.Cells.Item(33, 6) 161

# ---- Sheet 2: 演出 ----
$ws = $wb.Worksheets.Item(2)
Set-NumCell $ws.Cells.Item(4, 6) 719
Set-NumCell $ws.Cells.Item(5, 6) 607
Set-NumCell $ws.Cells.Item(6, 6) 607
Set-NumCell $ws.Cells.Item(9, 6) 12
Set-NumCell $ws.Cells.Item(12, 6) 269
Set-NumCell $ws.Cells.Item(15, 6) 330
Set-NumCell $ws.Cells.Item(16, 6) 330
Set-NumCell $ws.Cells.Item(17, 6) 67
Set-NumCell $ws.Cells.Item(19, 6) 928
Set-NumCell $ws.Cells.Item(26, 6) 227
Set-NumCell $ws.Cells.Item(27, 6) 221

# ---- Sheet 3: 本地生活 ----
$ws = $wb.Worksheets.Item(3)
Set-TextCell $ws.Cells.Item(5, 2) "2023-12-06"
Set-TextCell $ws.Cells.Item(5, 3) "上海·「咒术回战  × animate cafe」"
Set-TextCell $ws.Cells.Item(5, 4) "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
Set-TextCell $ws.Cells.Item(5, 5) "2023.12.06 00:00-2024.02.27 23:59"
Set-NumCell $ws.Cells.Item(5, 6) 2260
Set-NumCell $ws.Cells.Item(5, 7) 30
Set-TextCell $ws.Cells.Item(5, 8) "https://show.bilibili.com/platform/detail.html?id=79292"
Set-TextCell $ws.Cells.Item(5, 9) "//i2.hdslb.com/bfs/openplatform/202401/LyD46Kty1705488020552.jpeg"
Set-TextCell $ws.Cells.Item(6, 2) "2024-01-06"
Set-TextCell $ws.Cells.Item(6, 3) "上海·罗小黑 x HAPPY ZOO主题Cafe"
Set-TextCell $ws.Cells.Item(6, 4) "南京东路340号 百联zx创趣场"
Set-TextCell $ws.Cells.Item(6, 5) "2024.01.06 00:00-02.29 23:59"
Set-NumCell $ws.Cells.Item(6, 6) 898
Set-NumCell $ws.Cells.Item(6, 7) 10
Set-TextCell $ws.Cells.Item(6, 8) "https://show.bilibili.com/platform/detail.html?id=80171"
Set-TextCell $ws.Cells.Item(6, 9) "//i2.hdslb.com/bfs/openplatform/202312/chPePM8d1703485388734.png"
Set-TextCell $ws.Cells.Item(7, 2) "2024-01-21"
Set-TextCell $ws.Cells.Item(7, 3) "上海·JOYPOLIS世嘉都市乐园"
Set-TextCell $ws.Cells.Item(7, 4) "中山北路3300号 上海月星环球港"
Set-TextCell $ws.Cells.Item(7, 5) "2024.01.21 10:00-02.20 21:00"
Set-NumCell $ws.Cells.Item(7, 6) 1
Set-NumCell $ws.Cells.Item(7, 7) 190
Set-TextCell $ws.Cells.Item(7, 8) "https://show.bilibili.com/platform/detail.html?id=81140"
Set-TextCell $ws.Cells.Item(7, 9) "//i1.hdslb.com/bfs/openplatform/202401/7Bq5nJNe1705653236022.jpeg"
Set-TextCell $ws.Cells.Item(8, 2) "2024-01-22"
Set-TextCell $ws.Cells.Item(8, 3) "上海·「新春特惠」世嘉都市乐园-JP国潮杂技嘉年华"
Set-TextCell $ws.Cells.Item(8, 4) "中山北路3300号环球港购物中心4楼 上海世嘉都市乐园"
Set-TextCell $ws.Cells.Item(8, 5) "2024.01.22 14:00-03.03 18:40"
Set-NumCell $ws.Cells.Item(8, 6) 2
Set-NumCell $ws.Cells.Item(8, 7) 49
Set-TextCell $ws.Cells.Item(8, 8) "https://show.bilibili.com/platform/detail.html?id=81210"
Set-TextCell $ws.Cells.Item(8, 9) "//i2.hdslb.com/bfs/openplatform/202401/sw2khwYM1706086166106.jpeg"
Set-TextCell $ws.Cells.Item(9, 2) "2024-01-27"
Set-TextCell $ws.Cells.Item(9, 3) "上海・明日方舟主题店·[SWEET ZONE甜蜜区域]"
Set-TextCell $ws.Cells.Item(9, 4) "南京东路830号第一百货商业中心B馆5楼(海底捞旁边) 第一百货商业中心"
Set-TextCell $ws.Cells.Item(9, 5) "2024.01.27 00:00-03.31 23:59"
Set-NumCell $ws.Cells.Item(9, 6) 1115
Set-NumCell $ws.Cells.Item(9, 7) 30
Set-TextCell $ws.Cells.Item(9, 8) "https://show.bilibili.com/platform/detail.html?id=81277"
Set-TextCell $ws.Cells.Item(9, 9) "//i0.hdslb.com/bfs/openplatform/202401/hp6D0Drt1705991831205.jpeg"
Set-TextCell $ws.Cells.Item(10, 2) "2024-02-01"
Set-TextCell $ws.Cells.Item(10, 3) "上海·次元波板糖×线条小狗MALTESE 主题快闪店"
Set-TextCell $ws.Cells.Item(10, 4) "西藏北路166静安大悦城北座6楼611号 次元波板糖"
Set-TextCell $ws.Cells.Item(10, 5) "2024.02.01 00:00-03.01 23:59"
Set-NumCell $ws.Cells.Item(10, 6) 241
Set-NumCell $ws.Cells.Item(10, 7) 30
Set-TextCell $ws.Cells.Item(10, 8) "https://show.bilibili.com/platform/detail.html?id=81345"
Set-TextCell $ws.Cells.Item(10, 9) "//i0.hdslb.com/bfs/openplatform/202401/Qbpful951706080847394.png"
Set-TextCell $ws.Cells.Item(11, 2) "2024-02-02"
Set-TextCell $ws.Cells.Item(11, 3) "上海·2024《永远的7日之都》x  萌果酱谷子咖啡"
Set-TextCell $ws.Cells.Item(11, 4) "南京东路340号百联ZX 萌果酱谷子咖啡（百联）"
Set-TextCell $ws.Cells.Item(11, 5) "2024.02.02 00:00-03.10 23:59"
Set-NumCell $ws.Cells.Item(11, 6) 72
Set-NumCell $ws.Cells.Item(11, 7) 30
Set-TextCell $ws.Cells.Item(11, 8) "https://show.bilibili.com/platform/detail.html?id=81357"
Set-TextCell $ws.Cells.Item(11, 9) "//i2.hdslb.com/bfs/openplatform/202401/5OYoWSGL1706087914805.jpeg"
Set-TextCell $ws.Cells.Item(12, 2) "2024-02-15"
Set-TextCell $ws.Cells.Item(12, 3) "上海·飘起来吧魔法泡泡-魔术表演"
Set-TextCell $ws.Cells.Item(12, 4) "曹杨路1888号 上海露边社·演艺空间"
Set-TextCell $ws.Cells.Item(12, 5) "2024.02.15 19:00-03.03 20:10"
Set-NumCell $ws.Cells.Item(12, 6) 3
Set-NumCell $ws.Cells.Item(12, 7) 88
Set-TextCell $ws.Cells.Item(12, 8) "https://show.bilibili.com/platform/detail.html?id=81524"
Set-TextCell $ws.Cells.Item(12, 9) "//i0.hdslb.com/bfs/openplatform/202401/tls18D0J1706599640356.png"
# Row 13 no longer exists after the 2023-12-01 蔚蓝档案 entry was removed
# and subsequent rows B:I shifted up by one (column A indices stay fixed).
$ws.Rows.Item(13).Clear()

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item(4)
Set-NumCell $ws.Cells.Item(3, 6) 2260
Set-NumCell $ws.Cells.Item(5, 6) 1712
Set-NumCell $ws.Cells.Item(8, 6) 898
Set-NumCell $ws.Cells.Item(9, 6) 1115
Set-NumCell $ws.Cells.Item(10, 6) 241
Set-NumCell $ws.Cells.Item(11, 6) 72
Set-NumCell $ws.Cells.Item(12, 6) 719
Set-NumCell $ws.Cells.Item(13, 6) 1148
Set-NumCell $ws.Cells.Item(15, 6) 142
Set-NumCell $ws.Cells.Item(16, 6) 1401
Set-NumCell $ws.Cells.Item(17, 6) 607
Set-NumCell $ws.Cells.Item(18, 6) 63
Set-NumCell $ws.Cells.Item(19, 6) 93
Set-NumCell $ws.Cells.Item(20, 6) 627
Set-NumCell $ws.Cells.Item(23, 6) 87
Set-NumCell $ws.Cells.Item(24, 6) 1310
Set-NumCell $ws.Cells.Item(25, 6) 462
Set-NumCell $ws.Cells.Item(26, 6) 465
Set-NumCell $ws.Cells.Item(27, 6) 25
Set-NumCell $ws.Cells.Item(28, 6) 685
Set-NumCell $ws.Cells.Item(29, 6) 2539
Set-NumCell $ws.Cells.Item(31, 6) 16
Set-NumCell $ws.Cells.Item(32, 6) 283
Set-NumCell $ws.Cells.Item(33, 6) 269
Set-NumCell $ws.Cells.Item(35, 6) 97
Set-NumCell $ws.Cells.Item(37, 6) 562
Set-NumCell $ws.Cells.Item(38, 6) 918
Set-NumCell $ws.Cells.Item(39, 6) 330
Set-NumCell $ws.Cells.Item(40, 6) 67
Set-NumCell $ws.Cells.Item(42, 6) 63
Set-NumCell $ws.Cells.Item(46, 6) 227
Set-NumCell $ws.Cells.Item(47, 6) 221
Set-NumCell $ws.Cells.Item(49, 6) 161

